$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 16.4133
$ws.Range("B3").Value = 6.499200000000005
$ws.Range("E3").Value = 16.6961
$ws.Range("B4").Value = 8.631899999999996
$ws.Range("E9").Value = 17.52490000000002
$ws.Range("A11").Value = -21.6677
$ws.Range("A12").Value = -21.56770000000001
$ws.Range("B14").Value = 6.725399999999994
$ws.Range("A15").Value = -21.62209999999999
$ws.Range("E15").Value = 15.89490000000001
$ws.Range("E19").Value = 16.29489999999999
$ws.Range("E20").Value = 15.8059
$ws.Range("E25").Value = 17.18570000000001
$ws.Range("B26").Value = 4.722800000000002
$ws.Range("A27").Value = -21.94969999999999
$ws.Range("E27").Value = 16.77289999999999
$ws.Range("A28").Value = -21.90199999999999
$ws.Range("E28").Value = 16.5716
$ws.Range("E30").Value = 15.64080000000001
$ws.Range("A31").Value = -21.59909999999999
$ws.Range("B31").Value = 4.9588
$ws.Range("A32").Value = -21.36259999999999
$ws.Range("E32").Value = 16.59349999999999
$ws.Range("B35").Value = 9.115700000000007
$ws.Range("A36").Value = -19.16230000000001
$ws.Range("B37").Value = 8.739400000000003
$ws.Range("A38").Value = -18.99620000000001
$ws.Range("B39").Value = 9.369600000000004
$ws.Range("B40").Value = 8.851200000000004
$ws.Range("E44").Value = 16.91019999999999
$ws.Range("B45").Value = 6.819000000000004
$ws.Range("A46").Value = -21.72019999999999
$ws.Range("E47").Value = 16.7178
$ws.Range("B52").Value = 4.981100000000003
$ws.Range("A54").Value = -21.78269999999998
$ws.Range("A55").Value = -22.2716
$ws.Range("A56").Value = -21.97249999999999
$ws.Range("B57").Value = 4.642799999999992
$ws.Range("E58").Value = 16.29510000000001
$ws.Range("E62").Value = 16.71299999999999
$ws.Range("A67").Value = -21.50189999999998
$ws.Range("A69").Value = -21.67439999999997
$ws.Range("A72").Value = -22.12540000000003
$ws.Range("A73").Value = -19.22420000000001
$ws.Range("E77").Value = 17.21730000000003
$ws.Range("E78").Value = 16.62690000000003
$ws.Range("B81").Value = 6.535600000000001
$ws.Range("A83").Value = -21.75389999999998
$ws.Range("B83").Value = 5.977300000000001
$ws.Range("E84").Value = 16.7649
$ws.Range("A86").Value = -22.034
$ws.Range("E89").Value = 17.25080000000001
$ws.Range("A91").Value = -21.52540000000001
$ws.Range("E91").Value = 17.78090000000002
$ws.Range("E92").Value = 17.84520000000002
$ws.Range("A93").Value = -21.3173
$ws.Range("E96").Value = 16.74649999999999
$ws.Range("A99").Value = -20.00989999999999
$ws.Range("B100").Value = 4.956999999999997
$ws.Range("B102").Value = 8.474200000000005
$ws.Range("E102").Value = 16.90939999999998
